$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns updated for each data row (B, C, D, E, F, I, K)
$cols = @(2, 3, 4, 5, 6, 9, 11)

# New values for rows 2..25, one row per array, matching $cols order
$data = @(
    @(10.15175217277646, 7.591137565417178, 6.636188441705966, 16.31504681110695, 36.77115165189981, 28.50331753538431, 11.09120877063381),
    @(9.947993642223656, 7.305578535318016, 6.626819512963815, 15.39923975154763, 36.2728947711697, 28.26890353977672, 10.93826670383469),
    @(9.825094547100498, 7.127311570613896, 6.620946187909651, 14.81373973783727, 35.97052100670501, 28.12803250588032, 10.84749309914045),
    @(9.775654383845565, 7.054060930337483, 6.618521889491576, 14.56958767599031, 35.84831636770156, 28.0714312419368, 10.81134135927841),
    @(9.767486197857735, 7.041865241695828, 6.618117457203462, 14.5287199321262, 35.8280891667651, 28.06208225300221, 10.80539053999715),
    @(9.824425065198925, 7.126325954982144, 6.620913618670713, 14.81046911098311, 35.9688686447764, 28.1272658575191, 10.8470020797726),
    @(10.08108961371358, 7.493360986721551, 6.632982797120166, 16.00423940336702, 36.59868708436403, 28.42187816365396, 11.03785534076227),
    @(10.59815332256745, 8.184654637037099, 6.655715116606892, 18.17758843761273, 37.85662320908615, 29.02226980716632, 11.43464489911193),
    @(10.98145360590129, 8.669085516849917, 6.67187837026762, 19.80365641746459, 38.78754424259932, 29.47493120445546, 11.73669980367531),
    @(11.15556761358679, 8.883309059300874, 6.679121392432121, 20.50222134586855, 39.21098043157223, 29.68286424505762, 11.87571481866998),
    @(11.22138549805414, 8.96347185735643, 6.681848880390409, 20.76087198558848, 39.37119322328431, 29.76184915611497, 11.92853102948708),
    @(11.20721678304511, 8.9462511108414, 6.681262139194137, 20.70542792865362, 39.33669644997222, 29.74482815217328, 11.91714926070684),
    @(11.16098524004244, 8.889923717012159, 6.679346088097253, 20.52361834643397, 39.22416462248182, 29.68935772435232, 11.88005687099318),
    @(11.13264979017343, 8.855294605755279, 6.678170475090098, 20.41148992871486, 39.15521455182704, 29.65541108466844, 11.85735772757763),
    @(10.97006292632457, 8.654955173816898, 6.671402855880348, 19.75717871070606, 38.75985988446554, 29.46137903951733, 11.72764231236636),
    @(10.87019915457245, 8.530422552459187, 6.667223388641403, 19.34526785555455, 38.51722005245396, 29.34283121471895, 11.64843719397127),
    @(10.81274099451748, 8.458219986730995, 6.664809112956462, 19.10447308472871, 38.37766329313567, 29.27483801834718, 11.60303528837519),
    @(10.7932858390271, 8.433677184982153, 6.663989886885213, 19.02227757645139, 38.33041624344028, 29.25185108178405, 11.5876913241608),
    @(10.88083239359307, 8.543739316717181, 6.667669368414131, 19.38951721170741, 38.54305006859664, 29.35543124451823, 11.65685311165553),
    @(11.17456830442784, 8.906495016207911, 6.679909289755733, 20.57717954842836, 39.25722256270556, 29.70564443693237, 11.89094752541226),
    @(11.36583591518843, 9.137958790976167, 6.687819903982395, 21.31912438952093, 39.72314191584195, 29.93593717192304, 12.04493178158797),
    @(11.26384241929139, 9.014958637815861, 6.683605832078857, 20.92625772293325, 39.47458940734337, 29.81291158756348, 11.9626754450661),
    @(10.87602524075133, 8.537720692802237, 6.667467776618668, 19.36952448323944, 38.53137249501634, 29.34973426441891, 11.6530478529569),
    @(10.45732393523934, 8.001364724351724, 6.649662225365259, 17.59255190232296, 37.51461691874086, 28.85765071692624, 11.43464489911193)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $excelRow = $startRow + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($excelRow, $cols[$j]).Value = $rowValues[$j]
    }
}
